$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "G" = 0.5429463333333333; "H" = 1.628839; "I" = 0.04659251079363984; "J" = 0.04659251079363985; "M" = 4.877755666666666; "N" = 14.633267; "O" = 0.09961167132870688; "P" = 0.09961167132870689; "Q" = 2.648359554112555; "R" = 23.835235987013; "S" = 0.00464115787155528; "T" = 0.004641157871555281 }
    3 = @{ "G" = 0.5429463333333333; "H" = 1.628839; "I" = 0.04659251079363984; "J" = 0.04659251079363985; "O" = 0.1360673938501395; "P" = 0.1360673938501395; "Q" = 3.617602010883667; "R" = 32.558418097953; "S" = 0.006339721516625069; "T" = 0.00633972151662507 }
    4 = @{ "G" = 0.5429463333333333; "H" = 1.628839; "I" = 0.04659251079363984; "J" = 0.04659251079363985; "K" = 3; "L" = 1; "M" = 0.246459; "N" = 0.739377; "O" = 0.00503309197542868; "P" = 0.00503309197542868; "Q" = 0.133814010367; "R" = 1.204326093303; "S" = 0.0002345043921905428; "T" = 0.0002345043921905429 }
    5 = @{ "G" = 0.5429463333333333; "H" = 1.628839; "I" = 0.04659251079363984; "J" = 0.04659251079363985; "M" = 37.01331466666667; "N" = 111.039944; "O" = 0.7558718368280999; "P" = 0.7558718368280999; "Q" = 20.09624348277955; "R" = 180.866191345016; "S" = 0.03521796671602162; "T" = 0.03521796671602163 }
    6 = @{ "G" = 0.5429463333333333; "H" = 1.628839; "I" = 0.04659251079363984; "J" = 0.04659251079363985; "M" = 0.167274; "N" = 0.501822; "O" = 0.00341600601762507; "P" = 0.00341600601762507; "Q" = 0.090820804962; "R" = 0.8173872446579999; "S" = 0.0001591602972473347; "T" = 0.0001591602972473348 }
    7 = @{ "I" = 0.8858267105024722; "J" = 0.8858267105024723; "M" = 4.877755666666666; "N" = 14.633267; "O" = 0.09961167132870688; "P" = 0.09961167132870689; "Q" = 50.35117430004566; "R" = 453.160568700411; "S" = 0.08823867914076183; "T" = 0.08823867914076186 }
    8 = @{ "I" = 0.8858267105024722; "J" = 0.8858267105024723; "O" = 0.1360673938501395; "P" = 0.1360673938501395; "S" = 0.1205321319009134; "T" = 0.1205321319009134 }
    9 = @{ "I" = 0.8858267105024722; "J" = 0.8858267105024723; "K" = 3; "L" = 1; "M" = 0.246459; "N" = 0.739377; "O" = 0.00503309197542868; "P" = 0.00503309197542868; "Q" = 2.544100384449; "R" = 22.896903460041; "S" = 0.004458447308250377; "T" = 0.004458447308250377 }
    10 = @{ "I" = 0.8858267105024722; "J" = 0.8858267105024723; "M" = 37.01331466666667; "N" = 111.039944; "O" = 0.7558718368280999; "P" = 0.7558718368280999; "Q" = 382.0740491245947; "R" = 3438.666442121352; "S" = 0.6695714627788971; "T" = 0.6695714627788972 }
    11 = @{ "I" = 0.8858267105024722; "J" = 0.8858267105024723; "M" = 0.167274; "N" = 0.501822; "O" = 0.00341600601762507; "P" = 0.00341600601762507; "Q" = 1.726704432414; "R" = 15.540339891726; "S" = 0.003025989373649466; "T" = 0.003025989373649466 }
    12 = @{ "G" = 0.7875243333333334; "H" = 2.362573; "I" = 0.06758077870388791; "J" = 0.06758077870388793; "M" = 4.877755666666666; "N" = 14.633267; "O" = 0.09961167132870688; "P" = 0.09961167132870689; "Q" = 3.841351279554555; "R" = 34.57216151599101; "S" = 0.006731834316389756; "T" = 0.006731834316389759 }
    13 = @{ "G" = 0.7875243333333334; "H" = 2.362573; "I" = 0.06758077870388791; "J" = 0.06758077870388793; "O" = 0.1360673938501395; "P" = 0.1360673938501395; "Q" = 5.247202968285667; "R" = 47.22482671457101; "S" = 0.009195540432601037; "T" = 0.00919554043260104 }
    14 = @{ "G" = 0.7875243333333334; "H" = 2.362573; "I" = 0.06758077870388791; "J" = 0.06758077870388793; "K" = 3; "L" = 1; "M" = 0.246459; "N" = 0.739377; "O" = 0.00503309197542868; "P" = 0.00503309197542868; "Q" = 0.194092459669; "R" = 1.746832137021; "S" = 0.0003401402749877597; "T" = 0.0003401402749877597 }
    15 = @{ "G" = 0.7875243333333334; "H" = 2.362573; "I" = 0.06758077870388791; "J" = 0.06758077870388793; "M" = 37.01331466666667; "N" = 111.039944; "O" = 0.7558718368280999; "P" = 0.7558718368280999; "Q" = 29.14888595732356; "R" = 262.339973615912; "S" = 0.05108240733318109; "T" = 0.0510824073331811 }
    16 = @{ "G" = 0.7875243333333334; "H" = 2.362573; "I" = 0.06758077870388791; "J" = 0.06758077870388793; "M" = 0.167274; "N" = 0.501822; "O" = 0.00341600601762507; "P" = 0.00341600601762507; "Q" = 0.131732345334; "R" = 1.185591108006; "S" = 0.0002308563467282693; "T" = 0.0002308563467282693 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $addr = "{0}{1}" -f $col, $row
        $ws.Range($addr).Value = $updates[$row][$col]
    }
}

Write-Output "Update complete"